# Add the new "Wild Seven" rulebook item to the checklist.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row (row 3). The order below matches the order the
# values were originally typed in, which controls the order new entries
# land in the shared-strings table.
$ws.Range("E3").Value = "wild7.jpg"
$ws.Range("D3").Value = "Kokusai Tsushinsha"
$ws.Range("B3").Value = "ワイルド7"
$ws.Range("C3").Value = "Wild Seven"
$ws.Range("F3").Value = "rulebook"
$ws.Range("A3").Value = 2003

# Widen columns C (english) and D (publisher) so the new text fits.
$ws.Columns.Item(3).ColumnWidth = 24.498697916666668
$ws.Columns.Item(4).ColumnWidth = 22.666666666666668

# Leave the selection where it ends up after tabbing past the new row.
$ws.Range("F4").Select() | Out-Null
